$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.488.53'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '3.770.08'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D5").Value = "'616.37"
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").Value = "'177.67"
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '3.769.38'
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -1.64%  '
$ws.Range("E10").Value = '  -3.63%  '
$ws.Range("D11").Value = "'6.71"
$ws.Range("E11").Value = '  +6.05%  '
$ws.Range("D12").Value = "'0.483"
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").Value = "'40.04"
$ws.Range("E13").Value = '  -2.44%  '
$ws.Range("E14").Value = '  -4.26%  '
$ws.Range("D15").Value = '4.398.87'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").Value = '3.768.76'
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").Value = '69.572.19'
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = "'7.53"
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("E19").Value = '  -3.44%  '
$ws.Range("D20").Value = "'509.53"
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("D21").Value = "'16.35"
$ws.Range("E21").Value = '  -3.11%  '
$ws.Range("D22").Value = "'9.41"
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("D23").Value = "'0.728"
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").Value = "'86.48"
$ws.Range("E25").Value = '  -1.72%  '
$ws.Range("E26").Value = '  -2.63%  '
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("E28").Value = '  -3.77%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").Value = "'2.96"
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("D32").Value = "'7.97"
$ws.Range("E32").Value = '  +2.56%  '
$ws.Range("D33").Value = "'30.84"
$ws.Range("E33").Value = '  -2.07%  '
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("D37").Value = "'6.14"
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("D38").Value = "'0.139"
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("D39").Value = "'0.340"
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("D40").Value = "'450.98"
$ws.Range("E40").Value = '  +7.56%  '
$ws.Range("D42").Value = "'49.89"
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("D43").Value = "'2.97"
$ws.Range("E43").Value = '  +4.29%  '
$ws.Range("D44").Value = "'44.67"
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").Value = "'8.56"
$ws.Range("E45").Value = '  -2.28%  '
$ws.Range("D46").Value = '2.958.37'
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("E47").Value = '  -1.46%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "'27.16"
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").Value = "'138.69"
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("D51").Value = "'2.46"
$ws.Range("E51").Value = '  -0.48%  '
